# Update cryptocurrency price values (column D) on the active worksheet
# to reflect the latest scraped prices, as produced by the scheduled
# "Updated symbol list" GitHub Actions workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new price text.
$updates = @{
    "D2"  = "246.06"
    "D3"  = "22.01"
    "D4"  = "5.441"
    "D5"  = "0.05771"
    "D7"  = "6.351"
    "D8"  = "0.8184"
    "D9"  = "1.025"
    "D10" = "0.1430"
    "D11" = "0.07296"
    "D12" = "0.03104"
    "D14" = "4.165"
    "D15" = "0.09395"
    "D16" = "0.001597"
    "D17" = "0.04823"
    "D18" = "0.0005847"
    "D19" = "0.006300"
    "D20" = "0.004129"
    "D21" = "0.0009985"
    "D23" = "3.742"
    "D24" = "2.193"
    "D26" = "0.1329"
    "D27" = "0.0003997"
    "D40" = "0.03884"
    "D41" = "0.006699"
    "D42" = "0.1070"
    "D43" = "0.002299"
    "D44" = "0.006684"
    "D45" = "0.00005605"
    "D47" = "0.3898"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Leading apostrophe forces Excel to store the value as literal text
    # (matching the workbook's existing inline-string/text cells) instead
    # of auto-converting the numeric-looking string to a number.
    $cell.Value = "'" + $updates[$ref]
    # Reset the style back to Normal so no stray "quote prefix" / text
    # number-format style lingers on the cell (keeps formatting identical
    # to the original, unstyled data cells).
    $cell.Style = "Normal"
}
